$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab_5a_Indikatoren")

# Row 21: "17 bis 18-Jährigen" -> "17 und 18-Jährigen"
$ws.Range("D21").Value = "Anteil der 17 und 18-Jährigen mit (angestrebter) Studienberechtigung"
$ws.Range("E21").Value = "XXXAnteil der 17 und 18-Jährigen mit (angestrebter) Studienberechtigung"
$ws.Range("F21").Value = "Anteil der 17 und 18-Jährigen mit (angestrebter) Studienberechtigung"
$ws.Range("G21").Value = "XXXAnteil der 17 und 18-Jährigen mit (angestrebter) Studienberechtigung"
$ws.Range("L21").Value = "Anteil der 17 und 18-Jährigen mit (angestrebter) Studienberechtigung"
$ws.Range("M21").Value = "XXXAnteil der 17 und 18-Jährigen mit (angestrebter) Studienberechtigung"

# Row 48: spell out BIP
$ws.Range("H48").Value = "Jährlich mindestens 3,5 % des Bruttoinlandsprodukts (BIP) bis 2025"

# Row 50: Schulabsolvent* -> Schulabsolvierend*
$ws.Range("F50").Value = "Ausländische Schulabsolvierende"
$ws.Range("H50").Value = "Erhöhung des Anteils der ausländischen Schulabsolvierenden mit mindestens Hauptschulabschluss und Angleichung an die Quote deutscher Schulabsolvierender bis 2030"
$ws.Range("L50").Value = "Ausländische Schulabsolvierende"

# Row 58: rewording
$ws.Range("F58").Value = "Anteil der Personen mit hohen Wohnkosten"
